# Add "upi_vendeu" column (D) with Sim/Não values for the SP rows (2-124)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell D1 - match style of the other header cells (A1:C1): bold, centered
$ws.Range("D1").Value = "upi_vendeu"
$ws.Cells.Item(1, 4).Font.Bold = $true
$ws.Cells.Item(1, 4).HorizontalAlignment = -4108

# Values for D2:D124 (only SP-origin rows have data in this dataset)
$values = @(
  "Sim","Não","Não","Sim","Sim","Não","Não","Sim","Sim","Sim",
  "Sim","Não","Não","Não","Não","Sim","Sim","Não","Não","Sim",
  "Sim","Não","Sim","Não","Não","Não","Não","Sim","Sim","Não",
  "Sim","Não","Não","Sim","Sim","Não","Não","Sim","Sim","Sim",
  "Sim","Não","Sim","Não","Sim","Não","Não","Não","Não","Não",
  "Não","Não","Não","Sim","Não","Sim","Não","Não","Não","Não",
  "Sim","Não","Não","Não","Não","Não","Não","Não","Não","Não",
  "Não","Não","Não","Não","Não","Não","Não","Não","Não","Não",
  "Não","Sim","Não","Não","Não","Não","Não","Não","Não","Não",
  "Sim","Não","Sim","Sim","Não","Não","Não","Sim","Não","Não",
  "Não","Não","Não","Não","Não","Não","Sim","Não","Não","Não",
  "Sim","Não","Sim","Não","Sim","Não","Não","Não","Não","Sim",
  "Não","Não","Não"
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 4).Value = $values[$i]
}
